# Apply cryptos price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.787.72"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.891.18"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'0.7938"
$ws.Range("E5").Value = "  -2.33%  "
$ws.Range("D6").Value = "'241.82"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("D8").Value = "'0.3163"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").Value = "'25.64"
$ws.Range("E9").Value = "  -3.33%  "
$ws.Range("D10").Value = "'0.07036"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "'0.08050"
$ws.Range("D12").Value = "'0.7659"
$ws.Range("E12").Value = "  +2.65%  "
$ws.Range("D13").Value = "1.883.45"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").Value = "'5.321"
$ws.Range("D15").Value = "'92.33"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "29.765.37"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").Value = "'13.83"
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("D18").Value = "'5.941"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").Value = "'243.18"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").Value = "'0.000007696"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'8.201"
$ws.Range("E21").Value = "  +17.98%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'0.9999"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "2.142.71"
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "'0.1627"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Value = "'9.318"
$ws.Range("E26").Value = "  +1.01%  "
$ws.Range("D27").Value = "'164.03"
$ws.Range("E27").Value = "  -2.62%  "
$ws.Range("D28").Value = "'18.67"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("D29").Value = "'2.057"
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("D30").Value = "'1.379"
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("D31").Value = "'1.536"
$ws.Range("E31").Value = "  +1.69%  "
$ws.Range("D32").Value = "'4.455"
$ws.Range("E32").Value = "  +3.59%  "
$ws.Range("D33").Value = "'0.05730"
$ws.Range("E33").Value = "  +3.68%  "
$ws.Range("D34").Value = "'4.083"
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").Value = "'1.266"
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").Value = "'0.7387"
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("D39").Value = "'0.01914"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").Value = "'2.776"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("D41").Value = "'0.4403"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").Value = "'72.39"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "'5.840"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("D44").Value = "'0.8407"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "1.026.38"
$ws.Range("E46").Value = "  +4.18%  "
$ws.Range("D47").Value = "'102.21"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("D48").Value = "'9.914"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("D49").Value = "'1.851"
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("D50").Value = "'7.458"
$ws.Range("D51").Value = "2.054.41"
$ws.Range("E51").Value = "  -0.21%  "
